$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append results row for security/blowfish benchmark
$ws.Range("A14").Value = "security/blowfish/runme_large.sh"
$ws.Range("B14").Value = 0.14
$ws.Range("C14").Value = 0.13
$ws.Range("D14").Value = 0

# Match the author's final cursor position after entering the row
[void]$ws.Range("A16").Select()
